$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("PARTICIPANTS_TEMPLATE")

# Update the help text under "Language" column (R2) with the expanded wording.
$ws1.Range("R2").Value = "Include ISO-2 language code if different from PI language or leave blank if ISO-2 code is unknown"

# Update the active selection on the main sheet from A3 to R4.
$ws1.Range("R4").Select()
